$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EVALUACION2")

# Update the team members' full names (B4 / B5 on EVALUACION2, mirrored via
# formulas elsewhere in the workbook).
$ws.Range("B4").Value = "MORA PALMA MATIAS FRANCISCO"
$ws.Range("B5").Value = "LEYTON CISTERNA SEBASTIAN ANDRES"

# Change indicator #3's achieved category from "Completamente logrado" to
# "Logrado", which ripples through the weighted score formulas below it.
$ws.Range("C15").Value = "Logrado"

# Keep the active selection in sync with where the edit was made.
$ws.Range("D5").Select()
